$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 32: update title (D) and link (E)
$ws.Range("D32").Value = "아파치 피닉스(Apache Phoenix): HBase에서 SQL을 이용한다?!"
$ws.Range("E32").Value = "https://dodonam.tistory.com/511"

# Row 42: update title (D) and link (E)
$ws.Range("D42").Value = "[Linux-CentOS]네트워크 파일시스템 구조"
$ws.Range("E42").Value = "https://kjk92.tistory.com/128"

# Row 45: update title (D) only
$ws.Range("D45").Value = "RNN Auto-Encoder (RAE)"

# Row 47: update title (D) only
$ws.Range("D47").Value = "Pseudo Labeling, TTA(Test Time Augmentation) 기법"
